# The workbook tracks weekly wholesale-market price observations for
# "Zapallo italiano" at "Mercado Mayorista Lo Valledor de Santiago".
# The edit inserts one new weekly observation row at row 297 (pushing the
# existing rows 297-367 down to 298-368), and fills the new row with the
# new observation's data. Everything else (column headers, other rows)
# is left untouched; Excel's row Insert naturally shifts down formatting
# and the used-range "dimension" along with the cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 297 - this shifts old rows 297:367 to 298:368
# and extends the sheet's used range (dimension) to A1:R368 automatically.
$ws.Rows("297:297").Insert()

# Populate the newly inserted row 297 with the new observation.
$ws.Range("A297").Value = 6
$ws.Range("B297").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 44644
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112032
$ws.Range("G297").Value = "Zapallo italiano"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 230
$ws.Range("K297").Value = 10000
$ws.Range("L297").Value = 10000
$ws.Range("M297").Value = 10000
$ws.Range("N297").Value = "`$/caja 60 unidades"
$ws.Range("O297").Value = "Región Metropolitana"
$ws.Range("P297").Value = 167
$ws.Range("Q297").Value = 60
$ws.Range("R297").Value = "Hortaliza"
